# SpMV.xlsx update: fill in the numPipes = 2 benchmark data (rows 7-9,
# columns R:T) that used to be placeholder "nil" values, add the Ratio
# formulas (matching the pattern already used for numPipes = 1/4/8), and
# highlight the last "Ratio" row of every numPipes group in red - exactly
# like the other groups' last rows already do implicitly via the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7 (speed = 0.01 / numCols = 663552) ---
$ws.Range("R7").Value = 0.100855
$ws.Range("S7").Value = 0.019831000000000001
$ws.Range("T7").Formula = "=S7/S4"

# --- Row 8 (speed = 0.1 / numCols = 6709248), shared formula group T8:T9 ---
$ws.Range("R8").Value = 0.078479999999999994
$ws.Range("S8").Value = 0.025484
$ws.Range("T8").Formula = "=S8/S5"

# --- Row 9 (speed = 1 / numCols = 67108864) ---
$ws.Range("R9").Value = 0.076520000000000005
$ws.Range("S9").Value = 0.026137000000000001
$ws.Range("T9").Formula = "=S9/S6"

# Highlight the final "Ratio" cell of each numPipes group (the row where
# numCols = 67108864) in red, same as the new T9 row.
$ws.Range("T6").Font.Color = 255
$ws.Range("T9").Font.Color = 255
$ws.Range("T12").Font.Color = 255
$ws.Range("T15").Font.Color = 255

# Move the active selection to T9, matching where the new data now lives.
$ws.Range("T9").Select()
